# Feeds_Speeds_Vs_Material.xlsx - "Fixed some labels" edit
#
# 1) 2Flute sheet: shorten the MRR column header from the long endmill
#    description to just "Endmill", and clarify the MRR values as a
#    rate (".../min") instead of a bare volume.
# 2) Drop the trailing placeholder MRR values in rows 76-82 (they have
#    no corresponding data in the other columns).
# 3) Make the 2Flute tab the active/selected tab (instead of 3Flute),
#    with the selection parked at D2.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- 2Flute: column header / MRR units -------------------------------
$ws1.Range("C1").Value = "Endmill"

for ($r = 2; $r -le 75; $r++) {
    $cell = $ws1.Cells.Item($r, 10)
    if ($cell.Value2 -ne $null) {
        $cell.Value = $cell.Value2.ToString() + "/min"
    }
}

# Trailing rows 76-82 no longer carry an MRR placeholder value.
$ws1.Range("J76:J82").Clear()

# Column widths now reflect the shorter header / longer unit strings.
$ws1.Columns.Item(3).ColumnWidth = 16
$ws1.Columns.Item(10).ColumnWidth = 12.7

# --- Tab / selection state --------------------------------------------
$ws1.Activate()
$ws1.Range("D2").Select()
